$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 165, shifting existing rows 165-173 down to 166-174
$ws.Rows("165:165").Insert()

# Populate the newly inserted row 165 with the new data record
$ws.Range("A165").Value = 3
$ws.Range("B165").Value = "Femacal de La Calera"
$ws.Range("C165").Value = "Coquimbo"
$ws.Range("D165").Value = 44509
$ws.Range("E165").Value = 5
$ws.Range("F165").Value = 100112001
$ws.Range("G165").Value = "Berenjena"
$ws.Range("H165").Value = "Sin especificar"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 80
$ws.Range("K165").Value = 8000
$ws.Range("L165").Value = 8500
$ws.Range("M165").Value = 8250
$ws.Range("N165").Value = "$/caja 60 unidades"
$ws.Range("O165").Value = "Región de Arica y Parinacota"
$ws.Range("P165").Value = 138
$ws.Range("Q165").Value = 60
$ws.Range("R165").Value = "Hortaliza"
